# New crime data collected - update weekly CompStat figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings) ---
# "Volume 30   Number  48" -> "...49"
$ws.Range("A8").Value = "Volume 30   Number  49"

# "Report Covering the Week  11/27/2023  Through  12/3/2023"
#  -> "...12/4/2023  Through  12/10/2023"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Cells that change numeric type <-> text type ---
# G15: number 5 -> text "0"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G15").PasteSpecial(-4122)

# H15: number -100 -> text "***.*"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("N22").Copy()
$ws.Range("H15").PasteSpecial(-4122)

# C22: text "0" -> number 2
$ws.Range("J15").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

# C23: number 1 -> text "0"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)

# D28: number 1 -> text "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial(-4122)

# E28: number -100 -> text "***.*"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("N22").Copy()
$ws.Range("E28").PasteSpecial(-4122)

# D29: number 1 -> text "0"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial(-4122)

# E29: number -100 -> text "***.*"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("N22").Copy()
$ws.Range("E29").PasteSpecial(-4122)

# --- Plain numeric value updates ---
$ws.Range("G14").Value = 1

$ws.Range("N15").Value = 10

$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 75
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 12
$ws.Range("I16").Value = 310
$ws.Range("J16").Value = 291
$ws.Range("K16").Value = 6.529209621993
$ws.Range("L16").Value = 123.021582733813
$ws.Range("M16").Value = 57.360406091370
$ws.Range("N16").Value = -72.949389179755

$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -16.666666666666
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 4.347826086956
$ws.Range("I17").Value = 419
$ws.Range("J17").Value = 322
$ws.Range("K17").Value = 30.124223602484
$ws.Range("L17").Value = 115.979381443299
$ws.Range("M17").Value = 101.442307692308
$ws.Range("N17").Value = 26.969696969697

$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 21
$ws.Range("E18").Value = -52.380952380952
$ws.Range("F18").Value = 43
$ws.Range("G18").Value = 71
$ws.Range("H18").Value = -39.436619718309
$ws.Range("I18").Value = 529
$ws.Range("J18").Value = 560
$ws.Range("K18").Value = -5.535714285714
$ws.Range("L18").Value = 49.858356940509
$ws.Range("M18").Value = 15.250544662309
$ws.Range("N18").Value = -76.716549295774

$ws.Range("C19").Value = 31
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 63.157894736842
$ws.Range("F19").Value = 96
$ws.Range("G19").Value = 111
$ws.Range("H19").Value = -13.513513513513
$ws.Range("I19").Value = 1280
$ws.Range("J19").Value = 1417
$ws.Range("K19").Value = -9.668313338038
$ws.Range("L19").Value = 32.094943240454
$ws.Range("M19").Value = 98.757763975155
$ws.Range("N19").Value = -1.765157329240

$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 14
$ws.Range("E20").Value = -57.142857142857
$ws.Range("F20").Value = 32
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 485
$ws.Range("J20").Value = 274
$ws.Range("K20").Value = 77.007299270073
$ws.Range("L20").Value = 153.926701570681
$ws.Range("M20").Value = 73.214285714285
$ws.Range("N20").Value = -87.630706452435

$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = -7.8125
$ws.Range("F21").Value = 223
$ws.Range("G21").Value = 263
$ws.Range("H21").Value = -15.209125475285
$ws.Range("I21").Value = 3057
$ws.Range("J21").Value = 2910
$ws.Range("K21").Value = 5.051546391752
$ws.Range("L21").Value = 63.650963597430
$ws.Range("M21").Value = 68.801766979569
$ws.Range("N21").Value = -66.082325529790

$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 4
$ws.Range("I22").Value = 34
$ws.Range("J22").Value = 40
$ws.Range("K22").Value = -15
$ws.Range("L22").Value = 142.857142857143
$ws.Range("M22").Value = 580

$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -66.666666666666
$ws.Range("I23").Value = 17
$ws.Range("K23").Value = 6.25
$ws.Range("L23").Value = 21.428571428571
$ws.Range("M23").Value = 41.666666666666

$ws.Range("C24").Value = 53
$ws.Range("D24").Value = 60
$ws.Range("E24").Value = -11.666666666666
$ws.Range("F24").Value = 227
$ws.Range("G24").Value = 232
$ws.Range("H24").Value = -2.155172413793
$ws.Range("I24").Value = 2672
$ws.Range("J24").Value = 2639
$ws.Range("K24").Value = 1.250473664266
$ws.Range("L24").Value = 30.087633885102
$ws.Range("M24").Value = 78.490313961255

$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = -43.75
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 57
$ws.Range("H25").Value = 1.754385964912
$ws.Range("I25").Value = 842
$ws.Range("J25").Value = 673
$ws.Range("K25").Value = 25.111441307578
$ws.Range("L25").Value = 68.737474949899
$ws.Range("M25").Value = 23.823529411764

$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0

$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 94
$ws.Range("J27").Value = 79
$ws.Range("K27").Value = 18.987341772151
$ws.Range("L27").Value = 40.298507462686

$ws.Range("G28").Value = 2
$ws.Range("N28").Value = -69.230769230769

$ws.Range("G29").Value = 2
$ws.Range("N29").Value = -63.636363636363
